$d = $word.ActiveDocument

# --- Paragraph 3: " Check for instant search to expand when the word in the branch " ---
# -> " " (run1) + "Arrange form" (run2)
$p3 = $d.Paragraphs.Item(3)
$rng = $p3.Range
$rng.Find.Execute(" Check for instant search to expand when the word in the branch ", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2)
$p3 = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$insertPoint.InsertAfter("Arrange form")
# force a run split by toggling formatting on the newly inserted text, then
# resetting it back so the final formatting matches the rest of the run
$p3 = $d.Paragraphs.Item(3)
$newRunRng = $d.Range($p3.Range.End - 1 - 12, $p3.Range.End - 1)
$newRunRng.Font.Bold = 1
$newRunRng.Font.Bold = 0

Write-Output "Done paragraph 3"

# --- Paragraph 4: " Parse the query of SPARQL  and put it DIV" -> "Check for SKOS" ---
$p4 = $d.Paragraphs.Item(4)
$rng = $p4.Range
$rng.Find.Execute(" Parse the query of SPARQL  and put it DIV", $false, $false, $false, $false, $false, $true, 1, $false, "Check for SKOS", 2)

Write-Output "Done paragraph 4"

# --- Paragraph 5: " Graph representation of sparql result in neo4j" ---
# -> "Show graph in graphTab" (the " result in neo4j" trailing run is dropped,
# "sparql" is renamed to "graphTab" while keeping proofErr wrapper runs)
$p5 = $d.Paragraphs.Item(5)
$rng = $p5.Range
$rng.Find.Execute(" Graph representation of ", $false, $false, $false, $false, $false, $true, 1, $false, "Show graph in ", 2)
$p5 = $d.Paragraphs.Item(5)
$rng = $p5.Range
$rng.Find.Execute("sparql", $false, $false, $false, $false, $false, $true, 1, $false, "graphTab", 2)
$p5 = $d.Paragraphs.Item(5)
$rng = $p5.Range
$rng.Find.Execute(" result in neo4j", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "Done paragraph 5"

# --- Paragraph 6: "Check for the appearance (ba[bookmark]dge)" -> "Flow clean code style" ---
$p6 = $d.Paragraphs.Item(6)
$rng = $p6.Range
$rng.Find.Execute("Check for the appearance (badge)", $false, $false, $false, $false, $false, $true, 1, $false, "Flow clean code style", 2)

Write-Output "Done paragraph 6"

# --- New paragraph 7: "Click button on show file and show type" ---
$p6 = $d.Paragraphs.Item(6)
$newPara7 = $p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "Click button on show file and show type"

Write-Output "Done paragraph 7"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("[" + $i + "] " + $d.Paragraphs.Item($i).Range.Text)
}
